$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.00614754098360656,0.98155737704918,0.00204918032786885,0.00204918032786885,0.0163934426229508,0.00204918032786885,0.0122950819672131,0.973360655737705,0.0245901639344262,0.950819672131147,0.961065573770492,0,0.00409836065573771,1,0,0.995901639344262,0.00409836065573771,0.891393442622951,0.0614754098360656,0.069672131147541,0.989754098360656,0.00614754098360656,0.00204918032786885),
    @(0.0102459016393443,0.00409836065573771,0.0225409836065574,0.00614754098360656,0.00204918032786885,0,0.977459016393443,0.0102459016393443,0.014344262295082,0.00204918032786885,0.00409836065573771,0.944672131147541,0.00614754098360656,0,0.991803278688525,0.00409836065573771,0.993852459016393,0.0983606557377049,0.932377049180328,0.924180327868853,0.00204918032786885,0.00409836065573771,0.00409836065573771),
    @(0.979508196721312,0.00204918032786885,0.00819672131147541,0.00409836065573771,0.975409836065574,0.993852459016393,0.00819672131147541,0.00409836065573771,0.950819672131147,0.00409836065573771,0.0327868852459016,0.00204918032786885,0,0,0.00204918032786885,0,0.00204918032786885,0.00819672131147541,0.00204918032786885,0,0.00614754098360656,0.987704918032787,0.991803278688525),
    @(0.00409836065573771,0.0122950819672131,0.967213114754098,0.987704918032787,0.00614754098360656,0.00409836065573771,0.00204918032786885,0.0122950819672131,0.0102459016393443,0.040983606557377,0,0.0532786885245902,0.989754098360656,0,0.00614754098360656,0,0,0.00204918032786885,0.00409836065573771,0.00614754098360656,0.00204918032786885,0.00204918032786885,0.00204918032786885)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
